$wb = $excel.ActiveWorkbook

# Replace the "Score1".."Score5" text headers with plain numeric ratings 1-5
# on every sheet that has them. The column range differs per sheet depending
# on how many leading label columns it has.
$scoreSheets = @(
    @{ Name = "OIAbyband";             Range = "B1:F1" },
    @{ Name = "OIAbygenderethnicity "; Range = "C1:G1" },
    @{ Name = "OIAbyethnicity";        Range = "B1:F1" },
    @{ Name = "OIAbystatsnzethnicity"; Range = "C1:G1" },
    @{ Name = "OIAbygender";           Range = "B1:F1" }
)

foreach ($entry in $scoreSheets) {
    $ws = $wb.Worksheets.Item($entry.Name)
    $rng = $ws.Range($entry.Range)
    for ($i = 1; $i -le 5; $i++) {
        $rng.Cells.Item(1, $i).Value = $i
    }
}

# Update the saved selection on the sheet that was previously active...
$wsStats = $wb.Worksheets.Item("OIAbystatsnzethnicity")
$wsStats.Range("C1").Select()

# ...then switch the active tab to OIAbygender with its own saved selection.
$wsGender = $wb.Worksheets.Item("OIAbygender")
$wsGender.Activate()
$wsGender.Range("I23").Select()
